$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45758.01041666666, 2385),
    @(3, 45758.02083333334, 2384),
    @(4, 45758.03125, 2384),
    @(5, 45758.04166666666, 2379),
    @(6, 45758.05208333334, 2377),
    @(7, 45758.0625, 2379),
    @(8, 45758.07291666666, 2359),
    @(9, 45758.08333333334, 2357),
    @(10, 45758.09375, 2337),
    @(11, 45758.10416666666, 2335),
    @(12, 45758.11458333334, 2330),
    @(13, 45758.125, 2328),
    @(14, 45758.13541666666, 2203),
    @(15, 45758.14583333334, 2203),
    @(16, 45758.15625, 2199),
    @(17, 45758.16666666666, 2196),
    @(18, 45758.17708333334, 2134),
    @(19, 45758.1875, 2130),
    @(20, 45758.19791666666, 2139),
    @(21, 45758.20833333334, 2138),
    @(22, 45758.21875, 2027),
    @(23, 45758.22916666666, 2028),
    @(24, 45758.23958333334, 2030),
    @(25, 45758.25, 2027),
    @(26, 45758.26041666666, 1911),
    @(27, 45758.27083333334, 1923),
    @(28, 45758.28125, 1918),
    @(29, 45758.29166666666, 1917),
    @(30, 45758.30208333334, 1808),
    @(31, 45758.3125, 1808),
    @(32, 45758.32291666666, 1807),
    @(33, 45758.33333333334, 1808),
    @(34, 45758.34375, 1841),
    @(35, 45758.35416666666, 1847),
    @(36, 45758.36458333334, 1853),
    @(37, 45758.375, 1860),
    @(38, 45758.38541666666, 2083),
    @(39, 45758.39583333334, 2085),
    @(40, 45758.40625, 2086),
    @(41, 45758.41666666666, 2087),
    @(42, 45758.42708333334, 2174),
    @(43, 45758.4375, 2174),
    @(44, 45758.44791666666, 2173),
    @(45, 45758.45833333334, 2172),
    @(46, 45758.46875, 2209),
    @(47, 45758.47916666666, 2208),
    @(48, 45758.48958333334, 2207),
    @(49, 45758.5, 2206),
    @(50, 45758.51041666666, 1807),
    @(51, 45758.52083333334, 1806),
    @(52, 45758.53125, 1804),
    @(53, 45758.54166666666, 1803),
    @(54, 45758.55208333334, 1794),
    @(55, 45758.5625, 1792),
    @(56, 45758.57291666666, 1790),
    @(57, 45758.58333333334, 1788),
    @(58, 45758.59375, 1817),
    @(59, 45758.60416666666, 1816),
    @(60, 45758.61458333334, 1814),
    @(61, 45758.625, 1813),
    @(62, 45758.63541666666, 2061),
    @(63, 45758.64583333334, 2059),
    @(64, 45758.65625, 2056),
    @(65, 45758.66666666666, 2053),
    @(66, 45758.67708333334, 1969),
    @(67, 45758.6875, 1963),
    @(68, 45758.69791666666, 1957),
    @(69, 45758.70833333334, 1951),
    @(70, 45758.71875, 1596),
    @(71, 45758.72916666666, 1585),
    @(72, 45758.73958333334, 1573),
    @(73, 45758.75, 1561),
    @(74, 45758.76041666666, 1214),
    @(75, 45758.77083333334, 1206),
    @(76, 45758.78125, 1198),
    @(77, 45758.79166666666, 1191),
    @(78, 45758.80208333334, 940),
    @(79, 45758.8125, 935),
    @(80, 45758.82291666666, 931),
    @(81, 45758.83333333334, 926),
    @(82, 45758.84375, 731),
    @(83, 45758.85416666666, 729),
    @(84, 45758.86458333334, 726),
    @(85, 45758.875, 724),
    @(86, 45758.88541666666, 613),
    @(87, 45758.89583333334, 611),
    @(88, 45758.90625, 610),
    @(89, 45758.91666666666, 608),
    @(90, 45758.92708333334, 552),
    @(91, 45758.9375, 551),
    @(92, 45758.94791666666, 550),
    @(93, 45758.95833333334, 549),
    @(94, 45758.96875, 0),
    @(95, 45758.97916666666, 0),
    @(96, 45758.98958333334, 0),
    @(97, 45759.0, 0)
)

foreach ($item in $data) {
    $r = $item[0]
    $a = $item[1]
    $b = $item[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}
